$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename tfidf_* headers to tf-idf_*
$ws.Range("W1").Value = "tf-idf_mean"
$ws.Range("X1").Value = "tf-idf_std"
$ws.Range("Y1").Value = "tf-idf_fold0"
$ws.Range("Z1").Value = "tf-idf_fold1"
$ws.Range("AA1").Value = "tf-idf_fold2"
$ws.Range("AB1").Value = "tf-idf_fold3"
$ws.Range("AC1").Value = "tf-idf_fold4"
$ws.Range("AD1").Value = "tf-idf_chi_mean"
$ws.Range("AE1").Value = "tf-idf_chi_std"
$ws.Range("AF1").Value = "tf-idf_chi_fold0"
$ws.Range("AG1").Value = "tf-idf_chi_fold1"
$ws.Range("AH1").Value = "tf-idf_chi_fold2"
$ws.Range("AI1").Value = "tf-idf_chi_fold3"
$ws.Range("AJ1").Value = "tf-idf_chi_fold4"
$ws.Range("AK1").Value = "tf-idf_pca_mean"
$ws.Range("AL1").Value = "tf-idf_pca_std"
$ws.Range("AM1").Value = "tf-idf_pca_fold0"
$ws.Range("AN1").Value = "tf-idf_pca_fold1"
$ws.Range("AO1").Value = "tf-idf_pca_fold2"
$ws.Range("AP1").Value = "tf-idf_pca_fold3"
$ws.Range("AQ1").Value = "tf-idf_pca_fold4"

# Update recalculated metric/std values for kNN (row 2), RF (row 6), Ensemble (row 7)
$ws.Range("B2").Value = 0.6499308033204485
$ws.Range("C2").Value = 0.05188226497714451
$ws.Range("D2").Value = 0.5842397836666989
$ws.Range("E2").Value = 0.6034764308009091
$ws.Range("F2").Value = 0.6685715739762719
$ws.Range("G2").Value = 0.6630193883996419
$ws.Range("H2").Value = 0.7303468397587202
$ws.Range("P2").Value = 0.6868632838916314
$ws.Range("Q2").Value = 0.06004240274718479
$ws.Range("R2").Value = 0.6449637123403515
$ws.Range("S2").Value = 0.6945502170948722
$ws.Range("T2").Value = 0.6096117951153356
$ws.Range("U2").Value = 0.6977017551266377
$ws.Range("V2").Value = 0.7874889397809601
$ws.Range("W2").Value = 0.7890163061484834
$ws.Range("X2").Value = 0.03834817423520521
$ws.Range("Y2").Value = 0.8072686981715195
$ws.Range("Z2").Value = 0.7208008256395354
$ws.Range("AA2").Value = 0.8371644484547711
$ws.Range("AB2").Value = 0.7947466475692282
$ws.Range("AC2").Value = 0.7851009109073626
$ws.Range("AD2").Value = 0.7701386303190259
$ws.Range("AE2").Value = 0.0606711718265975
$ws.Range("AF2").Value = 0.7368890352140691
$ws.Range("AG2").Value = 0.6678966356385712
$ws.Range("AH2").Value = 0.8288865256607193
$ws.Range("AI2").Value = 0.8235085087061239
$ws.Range("AJ2").Value = 0.7935124463756458
$ws.Range("AK2").Value = 0.7747603453967136
$ws.Range("AL2").Value = 0.04669802744617439
$ws.Range("AM2").Value = 0.7747826891878354
$ws.Range("AN2").Value = 0.690707109042343
$ws.Range("AO2").Value = 0.7932054706248254
$ws.Range("AP2").Value = 0.8335110093312802
$ws.Range("AQ2").Value = 0.7815954487972837
$ws.Range("B6").Value = 0.824133966676615
$ws.Range("C6").Value = 0.05077201108534676
$ws.Range("D6").Value = 0.7786743270686259
$ws.Range("E6").Value = 0.7791195411464367
$ws.Range("F6").Value = 0.8405277892441185
$ws.Range("G6").Value = 0.8073995139463689
$ws.Range("H6").Value = 0.9149486619775244
$ws.Range("I6").Value = 0.8138681890195073
$ws.Range("J6").Value = 0.05308557299049709
$ws.Range("K6").Value = 0.7524983176757962
$ws.Range("L6").Value = 0.7773160676386482
$ws.Range("M6").Value = 0.8394693928328861
$ws.Range("N6").Value = 0.7965612957123993
$ws.Range("O6").Value = 0.9034958712378068
$ws.Range("P6").Value = 0.7534000713955357
$ws.Range("Q6").Value = 0.04603954068374024
$ws.Range("R6").Value = 0.6976906673655618
$ws.Range("S6").Value = 0.7235566106533849
$ws.Range("T6").Value = 0.7965620949047766
$ws.Range("U6").Value = 0.7307630000803294
$ws.Range("V6").Value = 0.8184279839736256
$ws.Range("W6").Value = 0.8051615125966769
$ws.Range("X6").Value = 0.02902032364832803
$ws.Range("Y6").Value = 0.7966558937312135
$ws.Range("Z6").Value = 0.777177128267368
$ws.Range("AA6").Value = 0.7851005362500172
$ws.Range("AB6").Value = 0.8074180170954365
$ws.Range("AC6").Value = 0.8594559876393493
$ws.Range("AD6").Value = 0.8206418104533046
$ws.Range("AE6").Value = 0.0150486189466485
$ws.Range("AF6").Value = 0.8292955480064367
$ws.Range("AG6").Value = 0.7994717977740049
$ws.Range("AH6").Value = 0.8296687449913256
$ws.Range("AI6").Value = 0.8062360201579217
$ws.Range("AJ6").Value = 0.8385369413368349
$ws.Range("AK6").Value = 0.8125355882640468
$ws.Range("AL6").Value = 0.03908838529513609
$ws.Range("AM6").Value = 0.7725007651058508
$ws.Range("AN6").Value = 0.7647826282759407
$ws.Range("AO6").Value = 0.8273405787045069
$ws.Range("AP6").Value = 0.8281230936578985
$ws.Range("AQ6").Value = 0.8699308755760369
$ws.Range("B7").Value = 0.8542573678287834
$ws.Range("C7").Value = 0.04194495864085813
$ws.Range("D7").Value = 0.8409703096153264
$ws.Range("G7").Value = 0.8062546772224193
$ws.Range("H7").Value = 0.9138696778350595
$ws.Range("I7").Value = 0.860939356328483
$ws.Range("J7").Value = 0.04645641239378697
$ws.Range("K7").Value = 0.8297463175122749
$ws.Range("M7").Value = 0.8925910075771353
$ws.Range("N7").Value = 0.8384982997886223
$ws.Range("P7").Value = 0.8282070169752302
$ws.Range("Q7").Value = 0.04355256613278493
$ws.Range("R7").Value = 0.7991180662430909
$ws.Range("S7").Value = 0.8109491118932504
$ws.Range("U7").Value = 0.8054788686446601
$ws.Range("W7").Value = 0.8773861715477101
$ws.Range("X7").Value = 0.0481244053602312
$ws.Range("Y7").Value = 0.8820442783208741
$ws.Range("AB7").Value = 0.85954598127773
$ws.Range("AC7").Value = 0.9132285358267801
$ws.Range("AD7").Value = 0.8599897504776187
$ws.Range("AE7").Value = 0.03728932460534133
$ws.Range("AG7").Value = 0.806017293310763
$ws.Range("AK7").Value = 0.8606441137027421
$ws.Range("AL7").Value = 0.0452830054076079
$ws.Range("AM7").Value = 0.8940327828905656
$ws.Range("AN7").Value = 0.7757105750314578
$ws.Range("AO7").Value = 0.9036182922886306
